# Add two new equipment items (rows 21 and 22) to the items sheet.
# The order in which text cells are assigned matters, since it determines
# the order new entries are appended to the shared string table.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 21: it_eq_rafishat ---
$ws.Range("Y21").Value = "res/assets/equipment/head/spritesheet_head_rafishat.png"
$ws.Range("D21").Value = "0, 250"
$ws.Range("A21").Value = "it_eq_rafishat"
$ws.Range("B21").Value = "itd_rafishat"
$ws.Range("C21").Value = 1
$ws.Range("E21").Value = 10
$ws.Range("H21").Value = -20
$ws.Range("I21").Value = 5
$ws.Range("O21").Value = 2
$ws.Range("P21").Value = 2
$ws.Range("Q21").Value = 1

# --- Row 22: it_eq_head_divinet3 ---
$ws.Range("Y22").Value = "res/assets/equipment/head/spritesheet_head_divinet3.png"
$ws.Range("A22").Value = "it_eq_head_divinet3"
$ws.Range("B22").Value = "itd_head_divinet3"
$ws.Range("D22").Value = "150, 200"
$ws.Range("C22").Value = 1
$ws.Range("E22").Value = 120
$ws.Range("F22").Value = 10
$ws.Range("G22").Value = 1
$ws.Range("H22").Value = 20
$ws.Range("I22").Value = 10
$ws.Range("K22").Value = 10
$ws.Range("N22").Value = 10
$ws.Range("O22").Value = 2
$ws.Range("S22").Value = 10

# Update selection to match the final workbook state
$ws.Range("S22").Select()
